# Update the Test 1..5 torque readings for the three Applied Torque rows
# (550, 350, 200) on the summary sheet. The source cells store these
# readings as literal text (e.g. "545.0"), not numbers, so the target
# range is pre-formatted as Text before the new values are written -
# otherwise Excel's automatic type inference would silently convert the
# numeric-looking strings into real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:G4").NumberFormat = "@"

# Row 2 (Applied Torque = 550)
$ws.Range("C2").Value = "545.0"
$ws.Range("D2").Value = "554.3"
$ws.Range("E2").Value = "558.7"
$ws.Range("F2").Value = "553.1"
$ws.Range("G2").Value = "560.8"

# Row 3 (Applied Torque = 350)
$ws.Range("C3").Value = "345.7"
$ws.Range("D3").Value = "350.2"
$ws.Range("E3").Value = "350.9"
$ws.Range("F3").Value = "354.2"
$ws.Range("G3").Value = "351.8"

# Row 4 (Applied Torque = 200)
$ws.Range("C4").Value = "197.7"
$ws.Range("D4").Value = "194.5"
$ws.Range("E4").Value = "194.0"
$ws.Range("F4").Value = "192.3"
$ws.Range("G4").Value = "192.3"
